# "Finish timing for now"
#
# The timing-diagram rounded-rectangle callout box (id=300,
# "Rounded Rectangle 299") on slide 1 is narrowed: its width shrinks
# from 627964 EMU to 470801 EMU while its height, position and all
# other formatting stay the same.
#
# PowerPoint's Shape sizing properties (.Left/.Top/.Width/.Height) are
# expressed in points, and OOXML stores them in EMU (914400 EMU per
# inch, 72 points per inch), so we convert explicitly to land on the
# exact target EMU value.

$EMUsPerPoint = 914400 / 72

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 300 -and $shp.Name -eq "Rounded Rectangle 299") {
        $target = $shp
        break
    }
}

if ($target -eq $null) {
    throw "Could not find shape 'Rounded Rectangle 299' (id=300) on slide 1"
}

$newWidthEmu = 470801
$target.Width = $newWidthEmu / $EMUsPerPoint
